$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 114
$ws.Range("A114").Value = 9
$ws.Range("B114").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C114").Value = "Metropolitana"
$ws.Range("D114").Value = 44628
$ws.Range("D114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E114").Value = 13
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100103
$ws.Range("H114").Value = "Frutos de hueso (carozo)"
$ws.Range("I114").Value = 100103002
$ws.Range("J114").Value = "Ciruela"
$ws.Range("K114").Value = "Angeleno"
$ws.Range("L114").Value = "Primera"
$ws.Range("M114").Value = 220
$ws.Range("N114").Value = 7000
$ws.Range("O114").Value = 7000
$ws.Range("P114").Value = 7000
$ws.Range("Q114").Value = "$/bandeja 18 kilos granel"
$ws.Range("R114").Value = "Región de O'Higgins"
$ws.Range("S114").Value = 389
$ws.Range("T114").Value = 18

# New row 115
$ws.Range("A115").Value = 9
$ws.Range("B115").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C115").Value = "Metropolitana"
$ws.Range("D115").Value = 44628
$ws.Range("D115").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E115").Value = 13
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100103
$ws.Range("H115").Value = "Frutos de hueso (carozo)"
$ws.Range("I115").Value = 100103002
$ws.Range("J115").Value = "Ciruela"
$ws.Range("K115").Value = "Angeleno"
$ws.Range("L115").Value = "Segunda"
$ws.Range("M115").Value = 300
$ws.Range("N115").Value = 5000
$ws.Range("O115").Value = 5000
$ws.Range("P115").Value = 5000
$ws.Range("Q115").Value = "$/bandeja 18 kilos granel"
$ws.Range("R115").Value = "Región de O'Higgins"
$ws.Range("S115").Value = 278
$ws.Range("T115").Value = 18
